# "change to another system" — the distance table is being trimmed down to
# drop the rows that route to the old system's nodes (INE, SA, MC) plus the
# now-redundant "Data Center -> FVT" and "PHM -> FVT" rows, leaving the
# remaining rows to shift up and fill the gaps, exactly like using Excel's
# own Delete Entire Row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete from the bottom up so the row numbers of the rows still to be
# removed don't shift out from under us.
$ws.Rows(31).Delete()   # PHM  -> FVT   (3, 10, 447.43)
$ws.Rows(16).Delete()   # Data Center -> MC   (1, 16, 693)
$ws.Rows(15).Delete()   # Data Center -> SA   (1, 15, 328)
$ws.Rows(14).Delete()   # Data Center -> INE  (1, 14, 474)
$ws.Rows(10).Delete()   # Data Center -> FVT  (1, 10, 447)

# Leave the workbook's selection/scroll roughly where the author left it.
$ws.Range("E25").Select()
